# Generate Report for Handback
# Updates the localization-status workbook:
#  - "Status" text changes from "Ready for handoff" to "Handback transform failed"
#    for the 644299c1-... row (Overview sheet + zh-cn/de-de detail sheets share the
#    same string, so editing the cells that use it updates all occurrences).
#  - Populates the "Error Detail" column (P) for that same row on the zh-cn and
#    de-de sheets with the handback/handoff file-name mismatch message.
#  - Widens the "Error Detail" column so the message is readable.

$wb = $excel.ActiveWorkbook

# --- Overview sheet: update the Status text shown for the 644299c1 row ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E3").Value = "Handback transform failed"
$wsOverview.Range("F3").Value = "Handback transform failed"

# Excel's ColumnWidth (characters) property has a constant ~5/6 character
# padding offset versus the raw OOXML <col width="..."> value, so to land on
# a stored width of exactly 40 we request 40 - 5/6.
$errorDetailColumnWidth = 40 - (5 / 6)

# --- zh-cn sheet ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C3").Value = "Handback transform failed"
$wsZhCn.Range("P3").Value = "Handback file name: ln1r1kfk.n4p is different with handoff file name: 644299c1-3979-4132-aa2e-8e1f9018b6fb.1a9464472e8cb5ba4de54f98b87ca88953887727.zh-cn."
$wsZhCn.Columns.Item(16).ColumnWidth = $errorDetailColumnWidth

# --- de-de sheet ---
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C3").Value = "Handback transform failed"
$wsDeDe.Range("P3").Value = "Handback file name: ln1r1kfk.n4p is different with handoff file name: 644299c1-3979-4132-aa2e-8e1f9018b6fb.1a9464472e8cb5ba4de54f98b87ca88953887727.de-de."
$wsDeDe.Columns.Item(16).ColumnWidth = $errorDetailColumnWidth

$wb.Save()
